$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row for the "Korpoelec Invasion" mission, marked as Advanced (column D)
$ws.Range("B9").Value = "Korpoelec Invasion"
$ws.Range("D9").Value = 1

$ws.Range("E13").Select()
